$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.937.52"
$ws.Range("D3").Value = "3.063.01"
$ws.Range("E3").Value = "  -0.78%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "'560.29"
$ws.Range("E5").Value = "  +0.20%  "
$ws.Range("D6").Value = "'142.65"
$ws.Range("E6").Value = "  -2.36%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").Value = "3.062.66"
$ws.Range("E8").Value = "  -0.67%  "
$ws.Range("E9").Value = "  +3.39%  "
$ws.Range("D10").Value = "'0.154"
$ws.Range("E10").Value = "  +0.47%  "
$ws.Range("D11").Value = "'6.12"
$ws.Range("E11").Value = "  -4.98%  "
$ws.Range("E12").Value = "  +1.88%  "
$ws.Range("E13").Value = "  +1.22%  "
$ws.Range("D14").Value = "'35.41"
$ws.Range("E14").Value = "  +0.27%  "
$ws.Range("D15").Value = "3.564.43"
$ws.Range("E15").Value = "  -0.63%  "
$ws.Range("D16").Value = "63.937.27"
$ws.Range("E16").Value = "  -0.93%  "
$ws.Range("D17").Value = "3.053.14"
$ws.Range("E17").Value = "  -1.08%  "
$ws.Range("E18").Value = "  -0.39%  "
$ws.Range("E19").Value = "  +0.00%  "
$ws.Range("D20").Value = "'487.92"
$ws.Range("E20").Value = "  +2.17%  "
$ws.Range("E21").Value = "  +2.79%  "
$ws.Range("D22").Value = "'0.693"
$ws.Range("E22").Value = "  +1.30%  "
$ws.Range("D23").Value = "'14.78"
$ws.Range("E23").Value = "  +8.91%  "
$ws.Range("D24").Value = "'7.53"
$ws.Range("E24").Value = "  -0.54%  "
$ws.Range("D25").Value = "'82.57"
$ws.Range("E25").Value = "  +1.55%  "
$ws.Range("D26").Value = "'1.00"
$ws.Range("E26").Value = "  -0.02%  "
$ws.Range("E27").Value = "  +0.64%  "
$ws.Range("D28").Value = "'8.21"
$ws.Range("E28").Value = "  -0.22%  "
$ws.Range("D29").Value = "'2.07"
$ws.Range("E29").Value = "  +0.11%  "
$ws.Range("D30").Value = "'1.00"
$ws.Range("E30").Value = "  +0.05%  "
$ws.Range("D31").Value = "'26.54"
$ws.Range("E31").Value = "  +1.35%  "
$ws.Range("E32").Value = "  -0.14%  "
$ws.Range("D33").Value = "'2.58"
$ws.Range("E33").Value = "  +3.51%  "
$ws.Range("D34").Value = "'5.78"
$ws.Range("E34").Value = "  +2.97%  "
$ws.Range("D35").Value = "'6.29"
$ws.Range("E35").Value = "  +1.29%  "
$ws.Range("D36").Value = "'54.90"
$ws.Range("E36").Value = "  -0.01%  "
$ws.Range("E37").Value = "  +1.19%  "
$ws.Range("D38").Value = "'442.89"
$ws.Range("E38").Value = "  -5.55%  "
$ws.Range("E39").Value = "  -2.22%  "
$ws.Range("D40").Value = "3.047.72"
$ws.Range("E40").Value = "  +2.73%  "
$ws.Range("E41").Value = "  +0.91%  "
$ws.Range("D42").Value = "'2.74"
$ws.Range("E42").Value = "  -8.76%  "
$ws.Range("E43").Value = "  +1.77%  "
$ws.Range("E44").Value = "  +6.40%  "
$ws.Range("D45").Value = "'28.02"
$ws.Range("E45").Value = "  -1.32%  "
$ws.Range("D46").Value = "'2.25"
$ws.Range("E46").Value = "  +3.81%  "
$ws.Range("E47").Value = "  -0.06%  "
$ws.Range("E48").Value = "  +1.17%  "
$ws.Range("D49").Value = "0.0₃0519"
$ws.Range("E49").Value = "  -0.70%  "
$ws.Range("D50").Value = "'117.43"
$ws.Range("E50").Value = "  -1.31%  "
$ws.Range("E51").Value = "  +2.97%  "
